$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.301
$ws.Range("D5").Value = 0.431
$ws.Range("E5").Value = 0.477
$ws.Range("F5").Value = 0.566
$ws.Range("G5").Value = 0.57
$ws.Range("H5").Value = 0.6

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.301
$ws.Range("E7").Value = 0.477
$ws.Range("F7").Value = 0.566
$ws.Range("H7").Value = 0.6

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.279
$ws.Range("D8").Value = 0.507
$ws.Range("E8").Value = 0.553
$ws.Range("F8").Value = 0.601
$ws.Range("G8").Value = 0.634
$ws.Range("H8").Value = 0.656

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.398
$ws.Range("C9").Value = 0.448
$ws.Range("D9").Value = 0.529
$ws.Range("E9").Value = 0.542
$ws.Range("F9").Value = 0.552
$ws.Range("G9").Value = 0.577
$ws.Range("H9").Value = 0.592
